# Auto-generated: apply cell-value updates from cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.562.12"
$ws.Range("E2").Value = "  +4.86%  "
$ws.Range("D3").Value = "2.488.69"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.77%  "
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "2.879.69"
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").Value = "2.486.96"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.847"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "47.444.88"
$ws.Range("E18").Value = "  +4.69%  "
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.77%  "
$ws.Range("D21").Value = "0.0₃0937"
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("E24").Value = "  +5.78%  "
$ws.Range("E25").Value = "  +3.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.32%  "
$ws.Range("E31").Value = "  +7.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("E34").Value = "  +2.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0783"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +5.31%  "
$ws.Range("E38").Value = "  +4.62%  "
$ws.Range("E39").Value = "  +4.71%  "
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.112"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.51%  "
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("D45").Value = "1.966.59"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.75%  "
